$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header J1 ("Residential" -> "Residentia") and add new column K ("houseex")
$ws.Range("J1").Value = "Residentia"
$ws.Range("K1").Value = "houseex"

# Default all K2:K35 to 0
for ($r = 2; $r -le 35; $r++) {
    $ws.Cells.Item($r, 11).Value = 0
}

# Specific non-zero "change in household expenditure" values for rows with data
$ws.Range("K7").Value = -4.9292343457456536
$ws.Range("K10").Value = -4.8203515127295331
$ws.Range("K14").Value = -7.5196914542577344
$ws.Range("K27").Value = -6.181305037989234
$ws.Range("K30").Value = -6.3433217044006769

# Update the view's selection to match the saved state in the file
$ws.Range("F15").Select()
